$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.594.80"
$ws.Range("E2").Value = "  -3.25%  "

$ws.Range("D3").Value = "2.607.24"
$ws.Range("E3").Value = "  -2.00%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.03"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.91%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.621"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.14%  "

$ws.Range("D9").Value = "2.604.59"
$ws.Range("E9").Value = "  -1.97%  "

$ws.Range("E10").Value = "  -7.93%  "

$ws.Range("E11").Value = "  -0.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.379"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -5.07%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.157"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.91"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.27%  "

$ws.Range("D15").Value = "3.074.59"
$ws.Range("E15").Value = "  -2.13%  "

$ws.Range("E16").Value = "  -7.75%  "

$ws.Range("D17").Value = "63.489.66"
$ws.Range("E17").Value = "  -3.27%  "

$ws.Range("D18").Value = "2.604.32"
$ws.Range("E18").Value = "  -2.57%  "

$ws.Range("E19").Value = "  -4.96%  "

$ws.Range("E20").Value = "  +0.22%  "

$ws.Range("E21").Value = "  -6.40%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "340.55"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.86%  "

$ws.Range("E23").Value = "  +0.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.05"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.38%  "

$ws.Range("E26").Value = "  -6.18%  "

$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.06"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.98%  "

$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "577.66"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.60%  "

$ws.Range("E29").Value = "  -4.45%  "

$ws.Range("E31").Value = "  -2.03%  "

$ws.Range("E32").Value = "  -3.84%  "

$ws.Range("E33").Value = "  -4.63%  "

$ws.Range("E34").Value = "  -5.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.32"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.24%  "

$ws.Range("E37").Value = "  -5.33%  "

$ws.Range("E38").Value = "  +0.03%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.60"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.78%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "153.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.14%  "

$ws.Range("E41").Value = "  -6.08%  "

$ws.Range("E42").Value = "  -0.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "41.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.24%  "

$ws.Range("E44").Value = "  -1.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "156.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.81%  "

$ws.Range("E46").Value = "  +0.37%  "

$ws.Range("E47").Value = "  -5.94%  "

$ws.Range("E48").Value = "  -6.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.628"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.60%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.96%  "

$ws.Range("E51").Value = "  -5.09%  "
